# Adding area ("Area"/"Atotal") calculations to the discharge worksheet,
# mirroring the existing Q / Qtotal layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new "Area" / "Atotal" headers -------------------
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Row 2: first area cell uses an explicit (non-shared) formula --------
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Row 3: also explicit, like the source E3 "first of range" cell ------
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# --- Rows 4-15: shared formula group (mirrors the D/E shared groups) -----
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Summary cells pulling the totals into J2 / K2 ------------------------
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Selection matches the authored state (J2:K2 active on J2) -----------
$ws.Range("J2:K2").Select()

$wb.Application.Calculate()
